$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the yearly table by one more column (N = year 2022), mirroring the
# existing column M formatting for rows 3-7 (header rule row, year header,
# and the three technology data rows).
$ws.Range("M3:M7").Copy($ws.Range("N3"))

$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 98.8
$ws.Range("N6").Value = 98
$ws.Range("N7").Value = 96.9

$excel.CutCopyMode = 0

# Leave the selection where it ended up after entering the new data.
$ws.Range("O4").Select()
